$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.451.43"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "3.689.87"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'235.92"
$ws.Range("E5").Value = "  -2.61%  "
$ws.Range("D6").Value = "'1.87"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("D7").Value = "'649.64"
$ws.Range("E7").Value = "  -0.94%  "
$ws.Range("D8").Value = "'0.427"
$ws.Range("E8").Value = "  +1.38%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  -2.56%  "
$ws.Range("D11").Value = "3.689.00"
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("D12").Value = "'0.0000308"
$ws.Range("E12").Value = "  +19.06%  "
$ws.Range("D13").Value = "'44.17"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "'6.73"
$ws.Range("E15").Value = "  +3.44%  "
$ws.Range("D16").Value = "4.379.21"
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").Value = "96.544.37"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "'8.80"
$ws.Range("E18").Value = "  +8.64%  "
$ws.Range("D19").Value = "3.729.65"
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("D20").Value = "'12.91"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "'18.65"
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("D22").Value = "'0.501"
$ws.Range("E22").Value = "  -5.94%  "
$ws.Range("D23").Value = "'517.58"
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("E24").Value = "  -1.80%  "
$ws.Range("D25").Value = "'0.0000209"
$ws.Range("E25").Value = "  +2.10%  "
$ws.Range("D26").Value = "'6.92"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").Value = "'100.58"
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("D28").Value = "'13.07"
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("D29").Value = "'0.177"
$ws.Range("E29").Value = "  +4.15%  "
$ws.Range("E30").Value = "  -0.77%  "
$ws.Range("D31").Value = "'12.06"
$ws.Range("E31").Value = "  +1.63%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").Value = "'1.85"
$ws.Range("E33").Value = "  +5.99%  "
$ws.Range("E34").Value = "  -1.04%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'653.70"
$ws.Range("E36").Value = "  +6.31%  "
$ws.Range("D37").Value = "'32.05"
$ws.Range("E37").Value = "  -3.25%  "
$ws.Range("D38").Value = "'0.586"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").Value = "'8.74"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  +12.67%  "
$ws.Range("D42").Value = "'2.08"
$ws.Range("E42").Value = "  +7.93%  "
$ws.Range("D43").Value = "'40.55"
$ws.Range("E43").Value = "  -3.76%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").Value = "'0.0448"
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("D47").Value = "'0.428"
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("D49").Value = "'2.27"
$ws.Range("E49").Value = "  -1.55%  "
$ws.Range("D50").Value = "'8.43"
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("D51").Value = "'3.51"
$ws.Range("E51").Value = "  +2.33%  "
